$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $text) {
    # Force the cell to keep a literal text value (e.g. "1", "11", a date
    # string, etc.) instead of having Excel auto-convert numeric-looking
    # text into a real number. Formatting as Text before assigning the
    # value, then clearing the format again afterwards, keeps the stored
    # value a string while leaving the cell's style untouched (style 0).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Set-EmptyPresentCell($cell) {
    # Touching the NumberFormat (and clearing it again) forces the cell to
    # actually be persisted in the sheet (and so counted in <dimension/>)
    # even though it carries no value - matching an empty-but-present
    # field in the source row.
    $cell.NumberFormat = "@"
    $cell.ClearFormats()
}

# --- Sheet "foresatt": add rows 2 and 3 ---
$wsF = $wb.Worksheets.Item("foresatt")

$wsF.Cells.Item(2, 1).Value = 0
$wsF.Cells.Item(2, 2).Value = 2
Set-TextValue $wsF.Cells.Item(2, 3) "1"
Set-TextValue $wsF.Cells.Item(2, 4) "1"
Set-TextValue $wsF.Cells.Item(2, 5) "1"
Set-TextValue $wsF.Cells.Item(2, 6) "1"

$wsF.Cells.Item(3, 1).Value = 1
$wsF.Cells.Item(3, 2).Value = 1
Set-TextValue $wsF.Cells.Item(3, 3) "1"
Set-TextValue $wsF.Cells.Item(3, 4) "1"
Set-TextValue $wsF.Cells.Item(3, 5) "1"
Set-TextValue $wsF.Cells.Item(3, 6) "1"

# copy the header-cell formatting (bold, border, centered) onto the new
# index cells in column A, matching the original workbook's style
$wsF.Cells.Item(1, 2).Copy()
$wsF.Range("A2:A3").PasteSpecial(-4122)

# --- Sheet "barn": add row 2 ---
$wsB = $wb.Worksheets.Item("barn")

$wsB.Cells.Item(2, 1).Value = 0
$wsB.Cells.Item(2, 2).Value = 1
Set-TextValue $wsB.Cells.Item(2, 3) "11"

$wsB.Cells.Item(1, 2).Copy()
$wsB.Range("A2").PasteSpecial(-4122)

# --- Sheet "soknad": add row 2 ---
$wsS = $wb.Worksheets.Item("soknad")

$wsS.Cells.Item(2, 1).Value = 0
$wsS.Cells.Item(2, 2).Value = 1
$wsS.Cells.Item(2, 3).Value = 2
$wsS.Cells.Item(2, 4).Value = 2
$wsS.Cells.Item(2, 5).Value = 1
$wsS.Cells.Item(2, 6).Value = "on"
Set-EmptyPresentCell $wsS.Cells.Item(2, 7)
Set-EmptyPresentCell $wsS.Cells.Item(2, 8)
Set-EmptyPresentCell $wsS.Cells.Item(2, 9)
Set-EmptyPresentCell $wsS.Cells.Item(2, 10)
Set-EmptyPresentCell $wsS.Cells.Item(2, 11)
Set-TextValue $wsS.Cells.Item(2, 12) "2024-11-06"
Set-TextValue $wsS.Cells.Item(2, 13) "1111"

$wsS.Cells.Item(1, 2).Copy()
$wsS.Range("A2").PasteSpecial(-4122)
